$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert this week's 4 new quality-grade rows at the top of the
# data block (row 18), pushing the previously-existing rows down by 4. The new
# rows start as a copy of what is now the most-recent prior week (now at rows
# 22:25) and then get this week's date + volumes written over them.
$ws.Rows.Item(18).Resize(4).Insert()

$ws.Range("A22:T25").Copy()
$ws.Range("A18").PasteSpecial()

# New reporting date for this week's rows.
$ws.Cells.Item(18,4).Value = 45126
$ws.Cells.Item(19,4).Value = 45126
$ws.Cells.Item(20,4).Value = 45126
$ws.Cells.Item(21,4).Value = 45126

# Updated volumes (column M) reported for this week.
$ws.Cells.Item(18,13).Value = 55
$ws.Cells.Item(19,13).Value = 60
$ws.Cells.Item(20,13).Value = 60
$ws.Cells.Item(21,13).Value = 45
